$d = $word.ActiveDocument

$replacements = @(
    @("431÷5=86, 1", "684÷2=342, 0"),
    @("311÷8=38, 7", "480÷4=120, 0"),
    @("239÷3=79, 2", "950÷8=118, 6"),
    @("373÷2=186, 1", "607÷3=202, 1"),
    @("699÷3=233, 0", "328÷5=65, 3"),
    @("154÷7=22, 0", "972÷2=486, 0"),
    @("634÷2=317, 0", "149÷9=16, 5"),
    @("680÷7=97, 1", "662÷8=82, 6"),
    @("380÷2=190, 0", "826÷8=103, 2"),
    @("272÷8=34, 0", "254÷3=84, 2"),
    @("470÷4=117, 2", "917÷4=229, 1"),
    @("430÷2=215, 0", "810÷5=162, 0"),
    @("576÷7=82, 2", "892÷8=111, 4"),
    @("949÷9=105, 4", "668÷9=74, 2"),
    @("487÷7=69, 4", "196÷4=49, 0"),
    @("472÷6=78, 4", "261÷2=130, 1"),
    @("746÷2=373, 0", "910÷9=101, 1"),
    @("168÷4=42, 0", "556÷7=79, 3"),
    @("774÷8=96, 6", "482÷6=80, 2"),
    @("572÷6=95, 2", "448÷3=149, 1"),
    @("997÷9=110, 7", "524÷2=262, 0"),
    @("750÷4=187, 2", "481÷9=53, 4"),
    @("503÷4=125, 3", "565÷4=141, 1"),
    @("408÷8=51, 0", "621÷5=124, 1"),
    @("652÷8=81, 4", "404÷9=44, 8")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
